# daily auto push: 2026-02-13 14:07 UTC
#
# A new daily-log entry for 2026/02/13 14:00 JST ("金" / 20 o'clock bucket,
# ranking 201) is inserted as a new data row right after the existing
# 2026/02/13 rows (old row 805), pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at 805; everything that used to live at row 805
# onward (old rows 805-846) shifts down to 806-847, preserving all of its
# data untouched.
$ws.Rows(805).Insert()

# Column A holds date-like text (e.g. "2026/02/13") that must stay a literal
# string, matching every other cell in the column, rather than being
# auto-converted into a real Excel date serial number. Temporarily force the
# cell to Text format while the value is entered, then clear the formatting
# back to the default so the cell ends up with no explicit style, exactly
# like its neighbours.
$newRow = $ws.Range("A805")
$newRow.NumberFormat = "@"
$newRow.Value = "2026/02/13"
$newRow.ClearFormats()

$ws.Range("B805").Value = "金"
$ws.Range("C805").Value = 20
$ws.Range("D805").Value = 201
